$d = $word.ActiveDocument

# 1) "Ruby is dynamic typed" -> "Ruby is strongly typed"
$d.Content.Find.Execute("Ruby is dynamic typed", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ruby is strongly typed", 2)

# 2) "string and character" -> "string and an integer “wat!”"
$d.Content.Find.Execute("string and character", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "string and an integer “wat!”", 2)
